$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1729
$ws.Range("J17").Value = 1625.1136
$ws.Range("L17").Value = 4875.3408
$ws.Range("N17").Value = -5211.3408

$ws.Range("H132").Value = 1582.1072
$ws.Range("I132").Value = 1550.0385
$ws.Range("K132").Value = 4650.1155
$ws.Range("M132").Value = -2120.1155

$ws.Range("H138").Value = 2472.7742
$ws.Range("I138").Value = 1798.5294
$ws.Range("J138").Value = 3291.5
$ws.Range("K138").Value = 5395.5882
$ws.Range("L138").Value = 9874.5
$ws.Range("M138").Value = -255.5882000000001
$ws.Range("N138").Value = -20154.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5743.706
$ws.Range("I2").Value = 6087
$ws.Range("J2").Value = 5438.5557
$ws.Range("K2").Value = 6087
$ws.Range("L2").Value = 5438.5557
$ws.Range("M2").Value = -5974
$ws.Range("N2").Value = -5664.5557

$ws.Range("H32").Value = 5595.5625
$ws.Range("I32").Value = 5595.5625
$ws.Range("K32").Value = 5595.5625
$ws.Range("M32").Value = -5308.5625

$ws.Range("H97").Value = 975.10345
$ws.Range("I97").Value = 1046.125
$ws.Range("J97").Value = 634.2
$ws.Range("K97").Value = 1046.125
$ws.Range("L97").Value = 634.2
$ws.Range("M97").Value = -550.125
$ws.Range("N97").Value = -1626.2

$ws.Range("H116").Value = 5743.706
$ws.Range("I116").Value = 6087
$ws.Range("J116").Value = 5438.5557
$ws.Range("K116").Value = 6087
$ws.Range("L116").Value = 5438.5557
$ws.Range("M116").Value = -3793
$ws.Range("N116").Value = -10026.5557

$ws.Range("H122").Value = 1833.1666
$ws.Range("I122").Value = 1833.1666
$ws.Range("K122").Value = 5499.4998
$ws.Range("M122").Value = -3049.4998

$ws.Range("H130").Value = 69713.5
$ws.Range("J130").Value = 69713.5
$ws.Range("L130").Value = 69713.5
$ws.Range("N130").Value = -79753.5

$ws.Range("H132").Value = 1584.2632
$ws.Range("I132").Value = 1390.6129
$ws.Range("K132").Value = 4171.8387
$ws.Range("M132").Value = -1641.8387

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5743.706
$ws.Range("I3").Value = 6087
$ws.Range("J3").Value = 5438.5557
$ws.Range("K3").Value = 6087
$ws.Range("L3").Value = 5438.5557
$ws.Range("M3").Value = -5973
$ws.Range("N3").Value = -5666.5557

$ws.Range("H20").Value = 1301.3448
$ws.Range("I20").Value = 1308.75
$ws.Range("J20").Value = 1284.8889
$ws.Range("K20").Value = 1308.75
$ws.Range("L20").Value = 1284.8889
$ws.Range("M20").Value = -1061.75
$ws.Range("N20").Value = -1778.8889

$ws.Range("H86").Value = 479476.25
$ws.Range("I86").Value = 2001464.8
$ws.Range("J86").Value = 3854.8125
$ws.Range("K86").Value = 2001464.8
$ws.Range("L86").Value = 3854.8125
$ws.Range("M86").Value = -2000341.8
$ws.Range("N86").Value = -6100.8125

$ws.Range("H89").Value = 479476.25
$ws.Range("I89").Value = 2001464.8
$ws.Range("J89").Value = 3854.8125
$ws.Range("K89").Value = 10007324
$ws.Range("L89").Value = 19274.0625
$ws.Range("M89").Value = -10001708
$ws.Range("N89").Value = -30506.0625

$ws.Range("H94").Value = 1145.3846
$ws.Range("I94").Value = 1190.8334
$ws.Range("K94").Value = 1190.8334
$ws.Range("M94").Value = -739.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2240.0527
$ws.Range("I16").Value = 2372.6875
$ws.Range("J16").Value = 1532.6666
$ws.Range("K16").Value = 2372.6875
$ws.Range("L16").Value = 1532.6666
$ws.Range("M16").Value = -2085.6875
$ws.Range("N16").Value = -2106.6666

$ws.Range("H31").Value = 3889.5652
$ws.Range("J31").Value = 3786.4666
$ws.Range("L31").Value = 3786.4666
$ws.Range("N31").Value = -4376.4666

$ws.Range("H34").Value = 3889.5652
$ws.Range("J34").Value = 3786.4666
$ws.Range("L34").Value = 3786.4666
$ws.Range("N34").Value = -4190.4666

$ws.Range("H113").Value = 2240.0527
$ws.Range("I113").Value = 2372.6875
$ws.Range("J113").Value = 1532.6666
$ws.Range("K113").Value = 2372.6875
$ws.Range("L113").Value = 1532.6666
$ws.Range("M113").Value = -202.6875
$ws.Range("N113").Value = -5872.6666

$ws.Range("H132").Value = 26484.812
$ws.Range("I132").Value = 19570.5
$ws.Range("K132").Value = 58711.5
$ws.Range("M132").Value = -56181.5

$ws.Range("H134").Value = 4757.5757
$ws.Range("I134").Value = 2815.111
$ws.Range("K134").Value = 8445.332999999999
$ws.Range("M134").Value = -5910.332999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 12500
$ws.Range("J106").Value = 12500
$ws.Range("L106").Value = 37500
$ws.Range("N106").Value = -39392

$ws.Range("H107").Value = 752.8
$ws.Range("J107").Value = 825.2857
$ws.Range("L107").Value = 2475.8571
$ws.Range("N107").Value = -6315.8571

$ws.Range("H113").Value = 1257.7
$ws.Range("I113").Value = 2966
$ws.Range("J113").Value = 830.625
$ws.Range("K113").Value = 8898
$ws.Range("L113").Value = 2491.875
$ws.Range("M113").Value = -6728
$ws.Range("N113").Value = -6831.875

$ws.Range("H117").Value = 1209.6666
$ws.Range("I117").Value = 1671.5
$ws.Range("J117").Value = 978.75
$ws.Range("K117").Value = 5014.5
$ws.Range("L117").Value = 2936.25
$ws.Range("M117").Value = -1572.5
$ws.Range("N117").Value = -9820.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6080.1
$ws.Range("I70").Value = 6760.8
$ws.Range("K70").Value = 6760.8
$ws.Range("M70").Value = -6490.8

$ws.Range("H73").Value = 6080.1
$ws.Range("I73").Value = 6760.8
$ws.Range("K73").Value = 6760.8
$ws.Range("M73").Value = -5824.8

$ws.Range("H102").Value = 1276
$ws.Range("I102").Value = 1276
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1276
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 346
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 3158.1428
$ws.Range("I122").Value = 3039.9092
$ws.Range("J122").Value = 3591.6667
$ws.Range("K122").Value = 9119.7276
$ws.Range("L122").Value = 10775.0001
$ws.Range("M122").Value = -6669.7276
$ws.Range("N122").Value = -15675.0001

$ws.Range("H132").Value = 14488.054
$ws.Range("I132").Value = 14006.167
$ws.Range("J132").Value = 15377.692
$ws.Range("K132").Value = 42018.501
$ws.Range("L132").Value = 46133.076
$ws.Range("M132").Value = -39488.501
$ws.Range("N132").Value = -51193.076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2200.4
$ws.Range("J22").Value = 2200.4
$ws.Range("L22").Value = 2200.4
$ws.Range("N22").Value = -2790.4

$ws.Range("H27").Value = 2200.4
$ws.Range("J27").Value = 2200.4
$ws.Range("L27").Value = 2200.4
$ws.Range("N27").Value = -2414.4

$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20540

$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -21872

$ws.Range("H132").Value = 7938.3335
$ws.Range("I132").Value = 8130.875
$ws.Range("J132").Value = 6398
$ws.Range("K132").Value = 24392.625
$ws.Range("L132").Value = 19194
$ws.Range("M132").Value = -21862.625
$ws.Range("N132").Value = -24254

$ws.Range("H136").Value = 5811.9443
$ws.Range("I136").Value = 5697.8335
$ws.Range("K136").Value = 17093.5005
$ws.Range("M136").Value = -14543.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4188.892
$ws.Range("I122").Value = 2548.3333
$ws.Range("J122").Value = 17723.5
$ws.Range("K122").Value = 7644.999899999999
$ws.Range("L122").Value = 53170.5
$ws.Range("M122").Value = -5194.999899999999
$ws.Range("N122").Value = -58070.5

$ws.Range("H132").Value = 137666.62
$ws.Range("I132").Value = 247850.19
$ws.Range("K132").Value = 743550.5700000001
$ws.Range("M132").Value = -741020.5700000001
